$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "H2-T23"
$ws.Cells.Item(2,3).Value = "Klrd1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 23.72136033333333
$ws.Cells.Item(2,8).Value = 71.164081
$ws.Cells.Item(2,9).Value = 0.2333880243108029
$ws.Cells.Item(2,10).Value = 0.2333880243108029
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.04040533333333333
$ws.Cells.Item(2,14).Value = 0.121216
$ws.Cells.Item(2,15).Value = 0.009158328487187175
$ws.Cells.Item(2,16).Value = 0.009158328487187173
$ws.Cells.Item(2,17).Value = 0.9584694713884445
$ws.Cells.Item(2,18).Value = 8.626225242496
$ws.Cells.Item(2,19).Value = 0.002137444191613959
$ws.Cells.Item(2,20).Value = 0.002137444191613958

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "H2-T23"
$ws.Cells.Item(3,3).Value = "Klrd1"
$ws.Cells.Item(3,4).Value = "M1"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 23.72136033333333
$ws.Cells.Item(3,8).Value = 71.164081
$ws.Cells.Item(3,9).Value = 0.2333880243108029
$ws.Cells.Item(3,10).Value = 0.2333880243108029
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 2.697694
$ws.Cells.Item(3,14).Value = 8.093081999999999
$ws.Cells.Item(3,15).Value = 0.6114630364782021
$ws.Cells.Item(3,16).Value = 0.611463036478202
$ws.Cells.Item(3,17).Value = 63.99297144307133
$ws.Cells.Item(3,18).Value = 575.9367429876419
$ws.Cells.Item(3,19).Value = 0.142708150022732
$ws.Cells.Item(3,20).Value = 0.142708150022732

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "H2-T23"
$ws.Cells.Item(4,3).Value = "Klrd1"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 23.72136033333333
$ws.Cells.Item(4,8).Value = 71.164081
$ws.Cells.Item(4,9).Value = 0.2333880243108029
$ws.Cells.Item(4,10).Value = 0.2333880243108029
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.673768333333333
$ws.Cells.Item(4,14).Value = 5.021305
$ws.Cells.Item(4,15).Value = 0.3793786350346109
$ws.Cells.Item(4,16).Value = 0.3793786350346108
$ws.Cells.Item(4,17).Value = 39.70406174952278
$ws.Cells.Item(4,18).Value = 357.336555745705
$ws.Cells.Item(4,19).Value = 0.08854243009645699
$ws.Cells.Item(4,20).Value = 0.08854243009645697

$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "H2-T23"
$ws.Cells.Item(5,3).Value = "Klrd1"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 13.94839233333333
$ws.Cells.Item(5,8).Value = 41.845177
$ws.Cells.Item(5,9).Value = 0.1372344453793459
$ws.Cells.Item(5,10).Value = 0.137234445379346
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.04040533333333333
$ws.Cells.Item(5,14).Value = 0.121216
$ws.Cells.Item(5,15).Value = 0.009158328487187175
$ws.Cells.Item(5,16).Value = 0.009158328487187173
$ws.Cells.Item(5,17).Value = 0.5635894416924444
$ws.Cells.Item(5,18).Value = 5.072304975232
$ws.Cells.Item(5,19).Value = 0.001256838130540996
$ws.Cells.Item(5,20).Value = 0.001256838130540996

$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "H2-T23"
$ws.Cells.Item(6,3).Value = "Klrd1"
$ws.Cells.Item(6,4).Value = "M1"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 13.94839233333333
$ws.Cells.Item(6,8).Value = 41.845177
$ws.Cells.Item(6,9).Value = 0.1372344453793459
$ws.Cells.Item(6,10).Value = 0.137234445379346
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 2.697694
$ws.Cells.Item(6,14).Value = 8.093081999999999
$ws.Cells.Item(6,15).Value = 0.6114630364782021
$ws.Cells.Item(6,16).Value = 0.611463036478202
$ws.Cells.Item(6,17).Value = 37.62849430727933
$ws.Cells.Item(6,18).Value = 338.656448765514
$ws.Cells.Item(6,19).Value = 0.08391379068105684
$ws.Cells.Item(6,20).Value = 0.08391379068105684

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "H2-T23"
$ws.Cells.Item(7,3).Value = "Klrd1"
$ws.Cells.Item(7,4).Value = "M2"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 13.94839233333333
$ws.Cells.Item(7,8).Value = 41.845177
$ws.Cells.Item(7,9).Value = 0.1372344453793459
$ws.Cells.Item(7,10).Value = 0.137234445379346
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.673768333333333
$ws.Cells.Item(7,14).Value = 5.021305
$ws.Cells.Item(7,15).Value = 0.3793786350346109
$ws.Cells.Item(7,16).Value = 0.3793786350346108
$ws.Cells.Item(7,17).Value = 23.34637738844278
$ws.Cells.Item(7,18).Value = 210.117396495985
$ws.Cells.Item(7,19).Value = 0.05206381656774812
$ws.Cells.Item(7,20).Value = 0.05206381656774813

$ws.Cells.Item(8,1).Value = "M1"
$ws.Cells.Item(8,2).Value = "H2-T23"
$ws.Cells.Item(8,3).Value = "Klrd1"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 23.93812733333334
$ws.Cells.Item(8,8).Value = 71.81438200000001
$ws.Cells.Item(8,9).Value = 0.2355207359746736
$ws.Cells.Item(8,10).Value = 0.2355207359746736
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 0.04040533333333333
$ws.Cells.Item(8,14).Value = 0.121216
$ws.Cells.Item(8,15).Value = 0.009158328487187175
$ws.Cells.Item(8,16).Value = 0.009158328487187173
$ws.Cells.Item(8,17).Value = 0.9672280142791113
$ws.Cells.Item(8,18).Value = 8.705052128512001
$ws.Cells.Item(8,19).Value = 0.002156976265600143
$ws.Cells.Item(8,20).Value = 0.002156976265600142

$ws.Cells.Item(9,1).Value = "M1"
$ws.Cells.Item(9,2).Value = "H2-T23"
$ws.Cells.Item(9,3).Value = "Klrd1"
$ws.Cells.Item(9,4).Value = "M1"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 23.93812733333334
$ws.Cells.Item(9,8).Value = 71.81438200000001
$ws.Cells.Item(9,9).Value = 0.2355207359746736
$ws.Cells.Item(9,10).Value = 0.2355207359746736
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 2.697694
$ws.Cells.Item(9,14).Value = 8.093081999999999
$ws.Cells.Item(9,15).Value = 0.6114630364782021
$ws.Cells.Item(9,16).Value = 0.611463036478202
$ws.Cells.Item(9,17).Value = 64.57774247836934
$ws.Cells.Item(9,18).Value = 581.199682305324
$ws.Cells.Item(9,19).Value = 0.1440122243726549
$ws.Cells.Item(9,20).Value = 0.1440122243726548

$ws.Cells.Item(10,1).Value = "M1"
$ws.Cells.Item(10,2).Value = "H2-T23"
$ws.Cells.Item(10,3).Value = "Klrd1"
$ws.Cells.Item(10,4).Value = "M2"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 23.93812733333334
$ws.Cells.Item(10,8).Value = 71.81438200000001
$ws.Cells.Item(10,9).Value = 0.2355207359746736
$ws.Cells.Item(10,10).Value = 0.2355207359746736
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 1.673768333333333
$ws.Cells.Item(10,14).Value = 5.021305
$ws.Cells.Item(10,15).Value = 0.3793786350346109
$ws.Cells.Item(10,16).Value = 0.3793786350346108
$ws.Cells.Item(10,17).Value = 40.06687948983445
$ws.Cells.Item(10,18).Value = 360.60191540851
$ws.Cells.Item(10,19).Value = 0.08935153533641865
$ws.Cells.Item(10,20).Value = 0.08935153533641864

$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "H2-T23"
$ws.Cells.Item(11,3).Value = "Klrd1"
$ws.Cells.Item(11,4).Value = "ECs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 38.75739366666667
$ws.Cells.Item(11,8).Value = 116.272181
$ws.Cells.Item(11,9).Value = 0.3813234742102279
$ws.Cells.Item(11,10).Value = 0.3813234742102279
$ws.Cells.Item(11,11).Value = 1
$ws.Cells.Item(11,12).Value = 0.3333333333333333
$ws.Cells.Item(11,13).Value = 0.04040533333333333
$ws.Cells.Item(11,14).Value = 0.121216
$ws.Cells.Item(11,15).Value = 0.009158328487187175
$ws.Cells.Item(11,16).Value = 0.009158328487187173
$ws.Cells.Item(11,17).Value = 1.566005410232889
$ws.Cells.Item(11,18).Value = 14.094048692096
$ws.Cells.Item(11,19).Value = 0.003492285636692715
$ws.Cells.Item(11,20).Value = 0.003492285636692714

$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "H2-T23"
$ws.Cells.Item(12,3).Value = "Klrd1"
$ws.Cells.Item(12,4).Value = "M1"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 38.75739366666667
$ws.Cells.Item(12,8).Value = 116.272181
$ws.Cells.Item(12,9).Value = 0.3813234742102279
$ws.Cells.Item(12,10).Value = 0.3813234742102279
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 2.697694
$ws.Cells.Item(12,14).Value = 8.093081999999999
$ws.Cells.Item(12,15).Value = 0.6114630364782021
$ws.Cells.Item(12,16).Value = 0.611463036478202
$ws.Cells.Item(12,17).Value = 104.5555883502047
$ws.Cells.Item(12,18).Value = 941.0002951518418
$ws.Cells.Item(12,19).Value = 0.2331652094210034
$ws.Cells.Item(12,20).Value = 0.2331652094210033

$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "H2-T23"
$ws.Cells.Item(13,3).Value = "Klrd1"
$ws.Cells.Item(13,4).Value = "M2"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 38.75739366666667
$ws.Cells.Item(13,8).Value = 116.272181
$ws.Cells.Item(13,9).Value = 0.3813234742102279
$ws.Cells.Item(13,10).Value = 0.3813234742102279
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 1.673768333333333
$ws.Cells.Item(13,14).Value = 5.021305
$ws.Cells.Item(13,15).Value = 0.3793786350346109
$ws.Cells.Item(13,16).Value = 0.3793786350346108
$ws.Cells.Item(13,17).Value = 64.87089820180056
$ws.Cells.Item(13,18).Value = 583.838083816205
$ws.Cells.Item(13,19).Value = 0.1446659791525319
$ws.Cells.Item(13,20).Value = 0.1446659791525319

$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "H2-T23"
$ws.Cells.Item(14,3).Value = "Klrd1"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 1.273876
$ws.Cells.Item(14,8).Value = 3.821628
$ws.Cells.Item(14,9).Value = 0.01253332012494962
$ws.Cells.Item(14,10).Value = 0.01253332012494962
$ws.Cells.Item(14,11).Value = 1
$ws.Cells.Item(14,12).Value = 0.3333333333333333
$ws.Cells.Item(14,13).Value = 0.04040533333333333
$ws.Cells.Item(14,14).Value = 0.121216
$ws.Cells.Item(14,15).Value = 0.009158328487187175
$ws.Cells.Item(14,16).Value = 0.009158328487187173
$ws.Cells.Item(14,17).Value = 0.05147138440533334
$ws.Cells.Item(14,18).Value = 0.463242459648
$ws.Cells.Item(14,19).Value = 0.0001147842627393624
$ws.Cells.Item(14,20).Value = 0.0001147842627393624

$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "H2-T23"
$ws.Cells.Item(15,3).Value = "Klrd1"
$ws.Cells.Item(15,4).Value = "M1"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 1.273876
$ws.Cells.Item(15,8).Value = 3.821628
$ws.Cells.Item(15,9).Value = 0.01253332012494962
$ws.Cells.Item(15,10).Value = 0.01253332012494962
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 2.697694
$ws.Cells.Item(15,14).Value = 8.093081999999999
$ws.Cells.Item(15,15).Value = 0.6114630364782021
$ws.Cells.Item(15,16).Value = 0.611463036478202
$ws.Cells.Item(15,17).Value = 3.436527641944
$ws.Cells.Item(15,18).Value = 30.928748777496
$ws.Cells.Item(15,19).Value = 0.007663661980755056
$ws.Cells.Item(15,20).Value = 0.007663661980755054

$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "H2-T23"
$ws.Cells.Item(16,3).Value = "Klrd1"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 1.273876
$ws.Cells.Item(16,8).Value = 3.821628
$ws.Cells.Item(16,9).Value = 0.01253332012494962
$ws.Cells.Item(16,10).Value = 0.01253332012494962
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 1.673768333333333
$ws.Cells.Item(16,14).Value = 5.021305
$ws.Cells.Item(16,15).Value = 0.3793786350346109
$ws.Cells.Item(16,16).Value = 0.3793786350346108
$ws.Cells.Item(16,17).Value = 2.132173309393333
$ws.Cells.Item(16,18).Value = 19.18955978454
$ws.Cells.Item(16,19).Value = 0.004754873881455206
$ws.Cells.Item(16,20).Value = 0.004754873881455205
